# Daily attendance processing - 2025-12-20 07:26:49
# Swap the order of the two comma-separated entries in the "Recorded By"
# column (G) whenever a cell contains exactly two values and the first
# value is not "backup@backdoor.com" (those triples/backup rows are left
# untouched, as are single-value cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ",\s*"

        if ($parts.Count -eq 2 -and $parts[0] -ne "backup@backdoor.com") {
            $newVal = $parts[1] + ", " + $parts[0]
            $cell.Value2 = $newVal
        }
    }
}
